# Append the 12-10-2025 gold-price row (row 21) to Sheet1, mirroring the
# existing data rows (dates in column A, price blurb in column B).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Leading apostrophe forces the date-like string to be stored as text
# instead of being auto-converted to a date serial number.
$ws.Range("A21").Value = "'12-10-2025"
$ws.Range("B21").Value = "The price of gold in India today is ₹12,508 per gram for 24 karat gold, ₹11,465 per gram for 22 karat gold and ₹9,381 per gram for 18 karat gold (also called 999 gold)."

# Copy the formatting (borders/fill/wrap) from the previous row so the new
# row matches the rest of the table, and so A21 ends up using the same
# style as the other date cells (rather than the quote-prefixed style that
# forcing text creates).
$ws.Range("A20:B20").Copy() | Out-Null
$ws.Range("A21:B21").PasteSpecial(-4122) | Out-Null
